$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) for rows 2-11 from 45208 to 45212
$ws.Range("C2:C11").Value = 45212
